$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 19:22"

# --- Estados Unidos (row 4): updated totals ---
$ws.Range("B4").Value = 619716
$ws.Range("C4").Value = 5830
$ws.Range("D4").Value = 47639
$ws.Range("E4").Value = 544887
$ws.Range("G4").Value = 1143
$ws.Range("H4").Value = 27190

# --- Turquia (row 12): updated totals ---
$ws.Range("B12").Value = 69392
$ws.Range("C12").Value = 4281
$ws.Range("D12").Value = 5674
$ws.Range("E12").Value = 62200
$ws.Range("F12").Value = 1820
$ws.Range("G12").Value = 115
$ws.Range("H12").Value = 1518

# --- Canada (row 15): updated totals ---
$ws.Range("B15").Value = 27593
$ws.Range("C15").Value = 530
$ws.Range("D15").Value = 8592
$ws.Range("E15").Value = 18047

# --- Suiza (row 16): updated totals ---
$ws.Range("E16").Value = 9697
$ws.Range("G16").Value = 65
$ws.Range("H16").Value = 1239

# --- Austria (row 20): updated totals (country stays the same) ---
$ws.Range("B20").Value = 14331
$ws.Range("C20").Value = 105
$ws.Range("E20").Value = 5840

# Rows 21-24 reshuffle: Irlanda jumps above India/Israel/Suecia, which
# each slide down one rank and keep their own prior totals.
# Row 21 becomes Irlanda with fresh totals.
$ws.Range("A21").Value = "Irlanda"
$ws.Range("B21").Value = 12547
$ws.Range("C21").Value = 1068
$ws.Range("D21").Value = 77
$ws.Range("E21").Value = 12026
$ws.Range("F21").Value = 194
$ws.Range("G21").Value = 38
$ws.Range("H21").Value = 444

# Row 22 becomes India (previously Israel), carrying old India's row-21 totals.
$ws.Range("A22").Value = "India"
$ws.Range("B22").Value = 12320
$ws.Range("C22").Value = 833
$ws.Range("D22").Value = 1432
$ws.Range("E22").Value = 10483
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = 405

# Row 23 becomes Israel (previously Suecia), carrying old Israel's row-22 totals.
$ws.Range("A23").Value = "Israel"
$ws.Range("B23").Value = 12200
$ws.Range("C23").Value = 154
$ws.Range("D23").Value = 2309
$ws.Range("E23").Value = 9765
$ws.Range("F23").Value = 176
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 126

# Row 24 becomes Suecia (previously Irlanda), carrying old Suecia's row-23 totals.
$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 11927
$ws.Range("C24").Value = 482
$ws.Range("D24").Value = 381
$ws.Range("E24").Value = 10343
$ws.Range("F24").Value = 954
$ws.Range("G24").Value = 170
$ws.Range("H24").Value = 1203

# --- Marruecos (row 60): updated totals ---
$ws.Range("B60").Value = 2024
$ws.Range("C60").Value = 136
$ws.Range("D60").Value = 229
$ws.Range("E60").Value = 1668

# Rows 102-104 reshuffle: Jordania jumps above Malta/Bolivia, which each
# slide down one rank and keep their own prior totals.
# Row 102 becomes Jordania with fresh totals.
$ws.Range("A102").Value = "Jordania"
$ws.Range("B102").Value = 401
$ws.Range("C102").Value = 4
$ws.Range("D102").Value = 250
$ws.Range("E102").Value = 144
$ws.Range("F102").Value = 5
$ws.Range("H102").Value = 7

# Row 103 becomes Malta (previously Bolivia), carrying old Malta's row-102 totals.
$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 399
$ws.Range("C103").Value = 6
$ws.Range("D103").Value = 44
$ws.Range("E103").Value = 352
$ws.Range("F103").Value = 4
$ws.Range("H103").Value = 3

# Row 104 becomes Bolivia (previously Jordania), carrying old Bolivia's row-103 totals.
$ws.Range("A104").Value = "Bolivia"
$ws.Range("C104").Value = 43
$ws.Range("D104").Value = 7
$ws.Range("E104").Value = 362
$ws.Range("F104").Value = 3
$ws.Range("H104").Value = 28
